$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.167.22'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.825.93'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('D4').Value = '''0.9997'
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').Value = '''241.63'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').Value = '''0.6181'
$ws.Range('E6').Value = '  -1.68%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').Value = '''0.07346'
$ws.Range('E8').Value = '  -2.42%  '
$ws.Range('D9').Value = '''0.2897'
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('D10').Value = '''22.99'
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('D11').Value = '''0.07675'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').Value = '1.813.43'
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('D13').Value = '''4.956'
$ws.Range('D14').Value = '''0.6621'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').Value = '''81.99'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '''0.000008907'
$ws.Range('E16').Value = '  -5.27%  '
$ws.Range('D17').Value = '''5.837'
$ws.Range('E17').Value = '  -2.52%  '
$ws.Range('D18').Value = '29.136.30'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').Value = '2.065.39'
$ws.Range('E19').Value = '  -0.81%  '
$ws.Range('D20').Value = '''237.24'
$ws.Range('E20').Value = '  +6.26%  '
$ws.Range('D21').Value = '''12.44'
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('D22').Value = '''1.0000'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '''7.151'
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('D24').Value = '''1.001'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').Value = '''0.1418'
$ws.Range('E26').Value = '  +1.71%  '
$ws.Range('D27').Value = '''8.436'
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('D28').Value = '''17.62'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('D29').Value = '''1.487'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').Value = '''0.05555'
$ws.Range('E30').Value = '  -3.85%  '
$ws.Range('D31').Value = '''4.093'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').Value = '''4.097'
$ws.Range('E32').Value = '  -1.31%  '
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('D34').Value = '''1.827'
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('D35').Value = '''0.7360'
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('D36').Value = '''1.133'
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('D38').Value = '''2.837'
$ws.Range('E38').Value = '  +2.86%  '
$ws.Range('D39').Value = '1.208.74'
$ws.Range('E39').Value = '  -1.33%  '
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = '''6.328'
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('D42').Value = '''0.9114'
$ws.Range('E42').Value = '  +2.68%  '
$ws.Range('D43').Value = '''0.9997'
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D44').Value = '''0.00000000130'
$ws.Range('E44').Value = '  +3.05%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '''101.42'
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.969.93'
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('D47').Value = '''64.71'
$ws.Range('E47').Value = '  -1.64%  '
$ws.Range('D48').Value = '''0.5085'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').Value = '''0.4009'
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('D50').Value = '''9.045'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('D51').Value = '''0.05761'
$ws.Range('E51').Value = '  -1.06%  '
